$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: C3:M3 were stored as text ("0", "200", ...); convert them to
# real numbers (matches row 2's style) ---
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

# --- New row 4: escalation entry, report totals for total count / defect
# count graphs ---
$ws.Range("A4").Value = "2025-02-24T12:54"
$ws.Range("B4").Value = "SHIFT_1"

# C4:M4 mirror row 3's original layout: numeric-looking values stored as
# TEXT. Force text typing with a leading apostrophe (like typing '200 in
# Excel), then reset the style so no extra formatting is left behind.
$ws.Range("C4").Value = "'200"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'1"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'0"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'0"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "'0"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = "'0"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "'0"
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Value = "'0"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "'0"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = "'0"
$ws.Range("M4").Style = "Normal"

$ws.Range("N4").Value = "suriya"
$ws.Range("O4").Value = "suriya"
$ws.Range("P4").Value = "suriya"
